# Auto-update: advance the "today" reference date by one day.
# Column E ("剩余" / days remaining) decrements by 1 for every data row,
# except the rows whose remaining count had already hit 1 - those shops
# were serviced today, so their cycle resets: column F ("开始时间") is
# set to the new "today" (2025-12-03) and column E goes back up to the
# full cycle length (column D, "总天").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose cycle resets today (remaining count was already 1).
$resetRows = @(50, 51, 52, 53, 54, 55, 56, 57)
$newStartDate = 20251203

foreach ($row in $resetRows) {
    $total = $ws.Cells.Item($row, 4).Value2   # column D - total days
    $ws.Cells.Item($row, 5).Value = $total    # column E - remaining resets to total
    $ws.Cells.Item($row, 6).Value = $newStartDate  # column F - new start date
}

# All other data rows (2-99) simply lose one day of remaining time,
# except row 36 whose start-date value is malformed/non-date and whose
# remaining already equals the total, so it is left untouched.
$skipRows = @(36) + $resetRows

for ($row = 2; $row -le 99; $row++) {
    if ($skipRows -contains $row) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)  # column E
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current - 1
    }
}
